$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 348-349, pushing the existing data
# (old rows 348..444) down to new rows 350..446.
$ws.Rows("348:349").Insert()

# Populate the two newly inserted rows with the new weekly records.
# Row 348: Calameño / Primera
$ws.Range("A348").Value = 4
$ws.Range("B348").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C348").Value = "Los Lagos"
$ws.Range("D348").Value = 44988
$ws.Range("E348").Value = 10
$ws.Range("F348").Value = 100112027
$ws.Range("G348").Value = "Melón"
$ws.Range("H348").Value = "Calameño"
$ws.Range("I348").Value = "Primera"
$ws.Range("J348").Value = 2000
$ws.Range("K348").Value = 1500
$ws.Range("L348").Value = 1500
$ws.Range("M348").Value = 1500
$ws.Range("N348").Value = "$/unidad"
$ws.Range("O348").Value = "Región de O'Higgins"
$ws.Range("P348").Value = 1500
$ws.Range("Q348").Value = 1
$ws.Range("R348").Value = "Hortaliza"

# Row 349: Tuna / Primera
$ws.Range("A349").Value = 4
$ws.Range("B349").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C349").Value = "Los Lagos"
$ws.Range("D349").Value = 44988
$ws.Range("E349").Value = 10
$ws.Range("F349").Value = 100112027
$ws.Range("G349").Value = "Melón"
$ws.Range("H349").Value = "Tuna"
$ws.Range("I349").Value = "Primera"
$ws.Range("J349").Value = 2000
$ws.Range("K349").Value = 1400
$ws.Range("L349").Value = 1400
$ws.Range("M349").Value = 1400
$ws.Range("N349").Value = "$/unidad"
$ws.Range("O349").Value = "Región de O'Higgins"
$ws.Range("P349").Value = 1400
$ws.Range("Q349").Value = 1
$ws.Range("R349").Value = "Hortaliza"
